$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("C2") "320018256721"
Set-TextValue $ws.Range("C3") "320018256732"
Set-TextValue $ws.Range("C4") "320018256765"
Set-TextValue $ws.Range("C5") "320018256787"
Set-TextValue $ws.Range("D5") "320018256787"
Set-TextValue $ws.Range("C6") "320018256824"
Set-TextValue $ws.Range("D6") "320018256824"
Set-TextValue $ws.Range("C7") "320018256846"
Set-TextValue $ws.Range("D7") "320018256846"
Set-TextValue $ws.Range("C8") "320018256879"
Set-TextValue $ws.Range("C9") "320018256890"
Set-TextValue $ws.Range("C10") "320018256927"
Set-TextValue $ws.Range("C11") "320018256949"
Set-TextValue $ws.Range("C12") "320018256982"
Set-TextValue $ws.Range("C13") "320018257007"
Set-TextValue $ws.Range("D13") "320018257007"
Set-TextValue $ws.Range("C14") "320018257030"
Set-TextValue $ws.Range("D14") "320018257030"
Set-TextValue $ws.Range("C15") "320018257051"
Set-TextValue $ws.Range("D15") "320018257051"
Set-TextValue $ws.Range("C16") "320018257084"
Set-TextValue $ws.Range("D16") "320018257084"
Set-TextValue $ws.Range("C17") "320018257100"
Set-TextValue $ws.Range("D17") "320018257100"
Set-TextValue $ws.Range("C18") "320018257143"
Set-TextValue $ws.Range("C19") "320018257165"
Set-TextValue $ws.Range("C20") "320018257198"
Set-TextValue $ws.Range("C21") "320018257213"
Set-TextValue $ws.Range("C22") "320018257246"

Write-Host "Edit applied"
